$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2388.2
$ws.Range("I138").Value = 1460.7037
$ws.Range("J138").Value = 2731.2466
$ws.Range("K138").Value = 4382.1111
$ws.Range("L138").Value = 8193.739799999999
$ws.Range("M138").Value = 757.8888999999999
$ws.Range("N138").Value = -18473.7398
$ws.Range("H139").Value = 69000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 69000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 69000
$ws.Range("N139").Value = -79280
$ws.Range("H141").Value = 7693.517
$ws.Range("I141").Value = 4396.857
$ws.Range("J141").Value = 100000
$ws.Range("K141").Value = 13190.571
$ws.Range("L141").Value = 300000
$ws.Range("M141").Value = -8010.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14378.384
$ws.Range("I32").Value = 14412.782
$ws.Range("J32").Value = 14000
$ws.Range("K32").Value = 14412.782
$ws.Range("L32").Value = 14000
$ws.Range("M32").Value = -14125.782
$ws.Range("H45").Value = 2999.6667
$ws.Range("I45").Value = 2999.6667
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2999.6667
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -2622.6667
$ws.Range("H61").Value = 2115.8572
$ws.Range("I61").Value = 2115.8572
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2115.8572
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1903.8572
$ws.Range("H74").Value = 1006
$ws.Range("I74").Value = 861
$ws.Range("J74").Value = 1206.7693
$ws.Range("K74").Value = 861
$ws.Range("L74").Value = 1206.7693
$ws.Range("M74").Value = 13
$ws.Range("N74").Value = -2954.7693
$ws.Range("H77").Value = 1006
$ws.Range("I77").Value = 861
$ws.Range("J77").Value = 1206.7693
$ws.Range("K77").Value = 4305
$ws.Range("L77").Value = 6033.8465
$ws.Range("M77").Value = 63
$ws.Range("N77").Value = -14769.8465
$ws.Range("H131").Value = 27250
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 27250
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 27250
$ws.Range("N131").Value = -37330
$ws.Range("H132").Value = 1872.3529
$ws.Range("I132").Value = 1522.0667
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 4566.2001
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -2036.2001
$ws.Range("N132").Value = -18558.5
$ws.Range("H133").Value = 40349
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 40349
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 40349
$ws.Range("N133").Value = -45409
$ws.Range("H136").Value = 2115.8572
$ws.Range("I136").Value = 2115.8572
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6347.571599999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3797.571599999999
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 71987.5
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 71987.5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 71987.5
$ws.Range("N132").Value = -82107.5
$ws.Range("H134").Value = 3625.75
$ws.Range("I134").Value = 3572.2856
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 10716.8568
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -8181.856800000001
$ws.Range("N134").Value = -17070
$ws.Range("H135").Value = 54724.285
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 54724.285
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 54724.285
$ws.Range("N135").Value = -64864.285
$ws.Range("H138").Value = 48000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 48000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 48000
$ws.Range("N138").Value = -58280
$ws.Range("H139").Value = 179887.73
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 179887.73
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 179887.73
$ws.Range("N139").Value = -190167.73
$ws.Range("H140").Value = 43500
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 43500
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 43500
$ws.Range("N140").Value = -53860
$ws.Range("H141").Value = 43593.8
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 43593.8
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 43593.8
$ws.Range("N141").Value = -53953.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 8907.23
$ws.Range("I4").Value = 2900
$ws.Range("J4").Value = 9999.454
$ws.Range("K4").Value = 2900
$ws.Range("L4").Value = 9999.454
$ws.Range("M4").Value = -2788
$ws.Range("N4").Value = -10223.454
$ws.Range("H31").Value = 28574658
$ws.Range("I31").Value = 41669172
$ws.Range("J31").Value = 4805
$ws.Range("K31").Value = 41669172
$ws.Range("L31").Value = 4805
$ws.Range("M31").Value = -41668877
$ws.Range("N31").Value = -5395
$ws.Range("H34").Value = 28574658
$ws.Range("I34").Value = 41669172
$ws.Range("J34").Value = 4805
$ws.Range("K34").Value = 41669172
$ws.Range("L34").Value = 4805
$ws.Range("M34").Value = -41668970
$ws.Range("N34").Value = -5209

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 977.5357
$ws.Range("I5").Value = 1587.375
$ws.Range("J5").Value = 733.6
$ws.Range("K5").Value = 4762.125
$ws.Range("L5").Value = 2200.8
$ws.Range("M5").Value = -4650.125
$ws.Range("N5").Value = -2424.8
$ws.Range("H132").Value = 1686.4445
$ws.Range("I132").Value = 994.1
$ws.Range("J132").Value = 2551.875
$ws.Range("K132").Value = 8946.9
$ws.Range("L132").Value = 22966.875
$ws.Range("M132").Value = -6416.9
$ws.Range("N132").Value = -28026.875
$ws.Range("H135").Value = 977.5357
$ws.Range("I135").Value = 1587.375
$ws.Range("J135").Value = 733.6
$ws.Range("K135").Value = 14286.375
$ws.Range("L135").Value = 6602.400000000001
$ws.Range("M135").Value = -11751.375
$ws.Range("N135").Value = -11672.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 18000
$ws.Range("I5").Value = 14000
$ws.Range("J5").Value = 30000
$ws.Range("K5").Value = 14000
$ws.Range("L5").Value = 30000
$ws.Range("M5").Value = -13888
$ws.Range("N5").Value = -30224
$ws.Range("H107").Value = 873.3333
$ws.Range("I107").Value = 1113.2667
$ws.Range("J107").Value = 273.5
$ws.Range("K107").Value = 1113.2667
$ws.Range("L107").Value = 273.5
$ws.Range("M107").Value = 806.7333000000001
$ws.Range("N107").Value = -4113.5
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2909.3635
$ws.Range("I68").Value = 1666.6666
$ws.Range("J68").Value = 4400.6
$ws.Range("K68").Value = 1666.6666
$ws.Range("L68").Value = 4400.6
$ws.Range("M68").Value = -917.6666
$ws.Range("N68").Value = -5898.6
$ws.Range("H71").Value = 2909.3635
$ws.Range("I71").Value = 1666.6666
$ws.Range("J71").Value = 4400.6
$ws.Range("K71").Value = 8333.333000000001
$ws.Range("L71").Value = 22003
$ws.Range("M71").Value = -4589.333000000001
$ws.Range("N71").Value = -29491
$ws.Range("H74").Value = 25000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 25000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 25000
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -26996
$ws.Range("H77").Value = 25000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 25000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 75000
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -84984
$ws.Range("H135").Value = 73429
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 73429
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 73429
$ws.Range("N135").Value = -83569

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 100002
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 100002
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 100002
$ws.Range("N5").Value = -100226
$ws.Range("H135").Value = 51579.918
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 51579.918
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 51579.918
$ws.Range("N135").Value = -61719.918
